# Revert "Powerpoint writer: consolidate text run nodes."
#
# Splits each word/space-delimited token of text ranges that currently
# hold multiple words glued together in a single <a:r> run back into
# separate runs -- one run per word and one run per inter-word space --
# matching the pre-consolidation run layout (e.g. "Slide " / "1 " /
# "(Content)" -> "Slide" / " " / "1" / " " / "(Content)").
#
# Mechanism: PowerPoint's TextRange.Characters(start, length) addresses a
# sub-range of an existing run. Assigning .Text on that sub-range (even to
# text identical to what is already there) forces the host to materialize
# that sub-range as its own run, splitting it off from whatever run used to
# span it. Walking the tokens left-to-right (skipping the final token, which
# is left as whatever remains of the trailing run) reproduces the desired
# per-word/per-space run boundaries without altering the visible text.
function Split-TextRunsBySpace($tr) {
    $full = $tr.Text
    $toks = [regex]::Matches($full, '[^ ]+| ')
    $offset = 1
    for ($i = 0; $i -lt $toks.Count - 1; $i++) {
        $tok = $toks[$i].Value
        $len = $toks[$i].Length
        $sub = $tr.Characters($offset, $len)
        $sub.Text = $tok
        $offset = $offset + $len
    }
}

$p = $ppt.ActivePresentation

# Slide titles: "Slide N (...)" style headings.
Split-TextRunsBySpace $p.Slides.Item(1).Shapes.Item("Title 1").TextFrame.TextRange
Split-TextRunsBySpace $p.Slides.Item(2).Shapes.Item("Title 1").TextFrame.TextRange
Split-TextRunsBySpace $p.Slides.Item(3).Shapes.Item("Title 1").TextFrame.TextRange
Split-TextRunsBySpace $p.Slides.Item(4).Shapes.Item("Title 1").TextFrame.TextRange
Split-TextRunsBySpace $p.Slides.Item(5).Shapes.Item("Title 1").TextFrame.TextRange
Split-TextRunsBySpace $p.Slides.Item(6).Shapes.Item("Title 1").TextFrame.TextRange
Split-TextRunsBySpace $p.Slides.Item(7).Shapes.Item("Title 1").TextFrame.TextRange
Split-TextRunsBySpace $p.Slides.Item(8).Shapes.Item("Title 1").TextFrame.TextRange
Split-TextRunsBySpace $p.Slides.Item(9).Shapes.Item("Title 1").TextFrame.TextRange
Split-TextRunsBySpace $p.Slides.Item(10).Shapes.Item("Title 1").TextFrame.TextRange
Split-TextRunsBySpace $p.Slides.Item(11).Shapes.Item("Title 1").TextFrame.TextRange
Split-TextRunsBySpace $p.Slides.Item(12).Shapes.Item("Title 1").TextFrame.TextRange

# Image captions ("an image" / "An image").
Split-TextRunsBySpace $p.Slides.Item(6).Shapes.Item("TextBox 3").TextFrame.TextRange
Split-TextRunsBySpace $p.Slides.Item(7).Shapes.Item("TextBox 3").TextFrame.TextRange
Split-TextRunsBySpace $p.Slides.Item(8).Shapes.Item("TextBox 3").TextFrame.TextRange
